$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Options")

# Delete the rows that correspond to references no longer needed:
# row 19 -> S874#1
# row 16 -> S837
# row 15 -> S801
# row 14 -> J802#
# row 11 -> H510#B
# row 10 -> H501#G1
# row 8  -> J572#65
# Delete from the bottom up so row numbers of the remaining rows don't shift
# before we get to them.
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(15).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(11).Delete()
$ws.Rows.Item(10).Delete()
$ws.Rows.Item(8).Delete()
